# "avances hasta mayo 2024"
# Applies the dataset update to the "data" sheet: adds a new `total_crt`
# column (AL), fixes the crt_1 column values that were stored as fractions
# (0.05 / 0.1) instead of the intended whole numbers (5 / 10), corrects two
# mis-typed second_dicotomic values (the literal letters "p"/"q" instead of
# their numeric codes), removes a stray numeric-format style from A9 and
# nudges the sheet view/selection to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$glosario = $wb.Worksheets.Item("glosario")

# --- New column header: AL1 = "total_crt" ---------------------------------
$ws.Range("AL1").Value = "total_crt"

# --- crt_1 (column F): stored as 0.05 / 0.1, should be 5 / 10 -------------
$ws.Range("F2").Value = 5
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 5
$ws.Range("F10").Value = 10
$ws.Range("F12").Value = 5
$ws.Range("F14").Value = 5

# --- second_dicotomic (column AK): fix mis-typed letter codes -------------
$ws.Range("AK14").Value = 5
$ws.Range("AK15").Value = 10

# --- New column values: total_crt (column AL) ------------------------------
$ws.Range("AL2").Value = 1
$ws.Range("AL3").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AL6").Value = 3
$ws.Range("AL7").Value = 0
$ws.Range("AL8").Value = 3
$ws.Range("AL9").Value = 0
$ws.Range("AL10").Value = 0
$ws.Range("AL11").Value = 1
$ws.Range("AL12").Value = 2
$ws.Range("AL13").Value = 0
$ws.Range("AL14").Value = 3
$ws.Range("AL15").Value = 0

# --- Drop the stray applyNumberFormat style that used to sit on A9 --------
$ws.Range("A9").ClearFormats()

# --- Sheet view / selection state ------------------------------------------
$ws.Activate()
$ws.Range("AO4").Select()

$glosario.Activate()
$glosario.Range("C22").Select()

$ws.Activate()
